# Work in progress: changes to Login, Stack and Linked List Module
#
# Adds six new worksheets (Sheet5..Sheet10) carrying QA test fixtures for a
# "Linked List" and "Stack" options page (dropdown option -> sample code ->
# expected console output / error), and moves the active tab off of Sheet2
# onto the newly added Sheet10.

$wb = $excel.ActiveWorkbook

function Add-SheetAfterLast {
    param([string]$name)
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
    $ws.Name = $name
    return $ws
}

function Fill-Rows {
    param($ws, $rows)
    for ($r = 0; $r -lt $rows.Count; $r++) {
        $row = $rows[$r]
        for ($c = 0; $c -lt $row.Count; $c++) {
            $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet5 / Sheet6 : blank placeholder sheets
# ---------------------------------------------------------------------------
$sheet5 = Add-SheetAfterLast "Sheet5"
$sheet6 = Add-SheetAfterLast "Sheet6"

# ---------------------------------------------------------------------------
# Sheet7 : OptionOnLinkedList -> code "abc123" -> NameError
# ---------------------------------------------------------------------------
$sheet7 = Add-SheetAfterLast "Sheet7"
$errAbc = "NameError: name 'abc123' is not defined on line 1"
Fill-Rows $sheet7 @(
    @("OptionOnLinkedList", "code", "errorMessage"),
    @("Introduction", "abc123", $errAbc),
    @("Creating Linked LIst", "abc123", $errAbc),
    @("Types of Linked List", "abc123", $errAbc),
    @("Implement Linked List in Python", "abc123", $errAbc),
    @("Traversal", "abc123", $errAbc),
    @("Insertion", "abc123", $errAbc),
    @("Deletion", "abc123", $errAbc)
)
[void]$sheet7.Range("C2").Select()

# ---------------------------------------------------------------------------
# Sheet8 : OptionOnLinkedList -> print('Hello World! <X>') -> Hello World! <X>
# ---------------------------------------------------------------------------
$sheet8 = Add-SheetAfterLast "Sheet8"
Fill-Rows $sheet8 @(
    @("OptionOnLinkedList", "code", "errorMessage"),
    @("Introduction", "print('Hello World! Introduction')", "Hello World! Introduction"),
    @("Creating Linked LIst", "print('Hello World! Creating Linked LIst')", "Hello World! Creating Linked LIst"),
    @("Types of Linked List", "print('Hello World! Types of Linked List')", "Hello World! Types of Linked List"),
    @("Implement Linked List in Python", "print('Hello World! Implement Linked List in Python')", "Hello World! Implement Linked List in Python"),
    @("Traversal", "print('Hello World! Traversal')", "Hello World! Traversal"),
    @("Insertion", "print('Hello World! Insertion')", "Hello World! Insertion"),
    @("Deletion", "print('Hello World! Deletion')", "Hello World! Deletion")
)
[void]$sheet8.Range("B4").Select()

# ---------------------------------------------------------------------------
# Sheet9 : OptionOnStack -> code "Abcd" -> NameError
# ---------------------------------------------------------------------------
$sheet9 = Add-SheetAfterLast "Sheet9"
$errAbcd = "NameError: name 'Abcd' is not defined on line 1"
Fill-Rows $sheet9 @(
    @("OptionOnStack", "code", "errorMessage"),
    @("Operations in Stack", "Abcd", $errAbcd),
    @("Implementation", "Abcd", $errAbcd),
    @("Applications", "Abcd", $errAbcd)
)
[void]$sheet9.Range("A4").Select()

# ---------------------------------------------------------------------------
# Sheet10 : OptionOnStack -> print('Hello World') -> Hello World
# ---------------------------------------------------------------------------
$sheet10 = Add-SheetAfterLast "Sheet10"
Fill-Rows $sheet10 @(
    @("OptionOnStack", "code", "message"),
    @("Operations in Stack", "print('Hello World')", "Hello World"),
    @("Implementation", "print('Hello World')", "Hello World"),
    @("Applications", "print('Hello World')", "Hello World")
)
[void]$sheet10.Range("C5").Select()

# Sheet10 ends up the active tab/sheet, matching the saved workbook view.
[void]$sheet10.Activate()
